# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New strikeout (K) counts per start, replacing the old "Strike#" (pitch-count of
# strikes thrown) values that used to live in column G.
$kValues = @{
    2  = 3
    3  = 1
    4  = 5
    5  = 5
    6  = 3
    7  = 0
    8  = 8
    9  = 5
    10 = 6
    11 = 4
    12 = 2
    13 = 4
    14 = 2
    15 = 0
    16 = 0
    17 = 0
    18 = 3
    19 = 2
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 2
    25 = 2
    26 = 2
    27 = 5
    28 = 1
    29 = 4
    30 = 1
    31 = 2
    32 = 5
    33 = 3
    34 = 1
    35 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
